# Updated symbol list on Mon Jan  9 07:44:49 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) quotes for the
# cryptos.xlsx coin tracker. Cells hold plain text (e.g. "278.00", "6.51%"),
# so force each cell to Text format before assigning the new string value -
# this stops Excel from re-interpreting the digits/percent sign as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.51%"
# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.76%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.810"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.64%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06247"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.45%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.930"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.79%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8799"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.52%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9414"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.81%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1451"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.41%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05258"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.56%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07345"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.68%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03112"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.47%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09062"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.08%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001566"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.40%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006278"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.63%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005861"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.72%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.452"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.20%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.266"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.10%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.264"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.56%"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3146"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.29%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1302"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.65%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.851"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.96%"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.32%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001175"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.21%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004283"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.06%"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001690"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.08%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04043"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.33%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006843"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "65.35%"
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.79%"
# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.46%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01211"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.59%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005086"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.47%"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.376"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "822.91%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
